$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation by prefixing values with a leading apostrophe,
# matching the original inline-string (text) cell contents from the source diff.

$ws.Range("D2").Value = "'29.366.55"
$ws.Range("D3").Value = "'1.840.84"
$ws.Range("E3").Value = "'  -0.17%  "
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "'  +0.15%  "
$ws.Range("D5").Value = "'238.97"
$ws.Range("E5").Value = "'  -0.42%  "
$ws.Range("D6").Value = "'0.6267"
$ws.Range("E6").Value = "'  -0.07%  "
$ws.Range("E7").Value = "'  +0.14%  "
$ws.Range("D8").Value = "'0.07420"
$ws.Range("E8").Value = "'  -0.82%  "
$ws.Range("D9").Value = "'0.2891"
$ws.Range("E9").Value = "'  -0.14%  "
$ws.Range("D10").Value = "'24.95"
$ws.Range("E10").Value = "'  +1.95%  "
$ws.Range("D11").Value = "'0.07723"
$ws.Range("D12").Value = "'1.836.39"
$ws.Range("E12").Value = "'  -0.43%  "
$ws.Range("D13").Value = "'4.970"
$ws.Range("E13").Value = "'  -0.26%  "
$ws.Range("D14").Value = "'0.6740"
$ws.Range("E14").Value = "'  -0.81%  "
$ws.Range("D15").Value = "'0.00001026"
$ws.Range("E15").Value = "'  -2.45%  "
$ws.Range("D16").Value = "'81.70"
$ws.Range("D17").Value = "'6.206"
$ws.Range("E17").Value = "'  +0.54%  "
$ws.Range("D18").Value = "'29.451.44"
$ws.Range("E18").Value = "'  +0.24%  "
$ws.Range("D19").Value = "'233.02"
$ws.Range("E19").Value = "'  +1.84%  "
$ws.Range("E20").Value = "'  +0.05%  "
$ws.Range("D22").Value = "'7.289"
$ws.Range("E22").Value = "'  -2.76%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "'  +0.18%  "
$ws.Range("D24").Value = "'157.82"
$ws.Range("E24").Value = "'  -0.41%  "
$ws.Range("D25").Value = "'8.499"
$ws.Range("E25").Value = "'  +0.88%  "
$ws.Range("E26").Value = "'  -1.71%  "
$ws.Range("E27").Value = "'  -1.11%  "
$ws.Range("D28").Value = "'0.07229"
$ws.Range("E28").Value = "'  +10.71%  "
$ws.Range("D29").Value = "'1.463"
$ws.Range("E29").Value = "'  +3.69%  "
$ws.Range("D30").Value = "'1.478"
$ws.Range("E30").Value = "'  -0.31%  "
$ws.Range("D31").Value = "'4.041"
$ws.Range("E31").Value = "'  -1.68%  "
$ws.Range("D32").Value = "'4.030"
$ws.Range("E32").Value = "'  -1.42%  "
$ws.Range("E33").Value = "'  -0.58%  "
$ws.Range("E34").Value = "'  -0.02%  "
$ws.Range("D35").Value = "'0.6965"
$ws.Range("E35").Value = "'  +0.54%  "
$ws.Range("D36").Value = "'2.569"
$ws.Range("E36").Value = "'  -0.50%  "
$ws.Range("E37").Value = "'  +0.42%  "
$ws.Range("D38").Value = "'6.907"
$ws.Range("E38").Value = "'  +1.92%  "
$ws.Range("D39").Value = "'2.817"
$ws.Range("E39").Value = "'  -0.71%  "
$ws.Range("D40").Value = "'1.235.29"
$ws.Range("E40").Value = "'  -2.05%  "
$ws.Range("D41").Value = "'0.9591"
$ws.Range("E41").Value = "'  +4.32%  "
$ws.Range("D42").Value = "'1.000"
$ws.Range("E42").Value = "'  +0.17%  "
$ws.Range("D43").Value = "'2.015.45"
$ws.Range("E43").Value = "'  +0.42%  "
$ws.Range("D44").Value = "'100.80"
$ws.Range("E44").Value = "'  -0.29%  "
$ws.Range("D45").Value = "'65.36"
$ws.Range("E45").Value = "'  -1.12%  "
$ws.Range("B46").Value = "'BabyDogeCoin"
$ws.Range("C46").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.00000000120"
$ws.Range("E46").Value = "'  +0.61%  "
$ws.Range("B47").Value = "'RenderToken"
$ws.Range("C47").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'1.714"
$ws.Range("E47").Value = "'  -0.61%  "
$ws.Range("B48").Value = "'Aptos"
$ws.Range("C48").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'6.936"
$ws.Range("E48").Value = "'  -1.90%  "
$ws.Range("B49").Value = "'EnergySwap"
$ws.Range("C49").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.848"
$ws.Range("E49").Value = "'  -1.61%  "
$ws.Range("B50").Value = "'TheSandbox"
$ws.Range("C50").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "'0.3899"
$ws.Range("E50").Value = "'  -1.32%  "
$ws.Range("B51").Value = "'Algorand"
$ws.Range("C51").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.1131"
$ws.Range("E51").Value = "'  -2.62%  "
